$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values are plain bracketed text (non-numeric) -- set directly
$textCells = @{
    "C2" = "[1.80; 3.42]"
    "G2" = "[1.35; 2.51]"
    "K2" = "[0.00; 0.00]"
    "C3" = "[1.04; 1.12]"
    "G3" = "[1.05; 1.12]"
    "K3" = "[1.05; 1.12]"
    "C4" = "[1.13; 1.91]"
    "G4" = "[1.11; 1.86]"
    "C5" = "[1.46; 9.09]"
    "G5" = "[1.62; 9.37]"
    "C6" = "[0.62; 2.10]"
    "C7" = "[0.45; 2.01]"
    "C8" = "[0.44; 4.85]"
    "C9" = "[0.62; 2.85]"
    "C10" = "[0.85; 0.92]"
    "G10" = "[0.84; 0.92]"
    "K10" = "[0.81; 0.89]"
}
foreach ($ref in $textCells.Keys) {
    $ws.Range($ref).Value = $textCells[$ref]
}

# Cells whose new values look like plain numbers but must stay stored as text
# (matching the source data which stores every value as a shared string).
# Force text storage via a temporary "@" (Text) number format, then restore the
# cell style to Normal so no stray per-cell style index remains.
$numericTextCells = @{
    "B2" = "2.48"
    "F2" = "1.84"
    "J2" = "0.00"
    "B3" = "1.08"
    "D3" = "20.6"
    "F3" = "1.08"
    "H3" = "22.5"
    "J3" = "1.08"
    "L3" = "28.3"
    "B4" = "1.47"
    "D4" = "8.1"
    "F4" = "1.44"
    "H4" = "7.8"
    "B5" = "3.64"
    "D5" = "7.7"
    "F5" = "3.90"
    "H5" = "9.3"
    "B6" = "1.14"
    "B7" = "0.95"
    "B8" = "1.47"
    "B9" = "1.33"
    "D9" = "0.5"
    "B10" = "0.88"
    "F10" = "0.88"
    "J10" = "0.85"
}
foreach ($ref in $numericTextCells.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $numericTextCells[$ref]
    $cell.Style = "Normal"
}
